$d = $word.ActiveDocument

$d.Content.Find.Execute("2023-09-01 Friday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-09-02 Saturday", 2) | Out-Null
$d.Content.Find.Execute("25×76=1900", $true, $false, $false, $false, $false, $true, 1, $false, "36×99=3564", 2) | Out-Null
$d.Content.Find.Execute("52×15=780", $true, $false, $false, $false, $false, $true, 1, $false, "61×97=5917", 2) | Out-Null
$d.Content.Find.Execute("84×25=2100", $true, $false, $false, $false, $false, $true, 1, $false, "82×31=2542", 2) | Out-Null
$d.Content.Find.Execute("44×28=1232", $true, $false, $false, $false, $false, $true, 1, $false, "13×94=1222", 2) | Out-Null
$d.Content.Find.Execute("39×40=1560", $true, $false, $false, $false, $false, $true, 1, $false, "50×68=3400", 2) | Out-Null
$d.Content.Find.Execute("62×73=4526", $true, $false, $false, $false, $false, $true, 1, $false, "71×85=6035", 2) | Out-Null
$d.Content.Find.Execute("52×16=832", $true, $false, $false, $false, $false, $true, 1, $false, "57×98=5586", 2) | Out-Null
$d.Content.Find.Execute("62×46=2852", $true, $false, $false, $false, $false, $true, 1, $false, "87×51=4437", 2) | Out-Null
$d.Content.Find.Execute("85×38=3230", $true, $false, $false, $false, $false, $true, 1, $false, "75×84=6300", 2) | Out-Null
$d.Content.Find.Execute("50×42=2100", $true, $false, $false, $false, $false, $true, 1, $false, "74×26=1924", 2) | Out-Null
$d.Content.Find.Execute("28×56=1568", $true, $false, $false, $false, $false, $true, 1, $false, "44×79=3476", 2) | Out-Null
$d.Content.Find.Execute("66×35=2310", $true, $false, $false, $false, $false, $true, 1, $false, "23×64=1472", 2) | Out-Null
$d.Content.Find.Execute("37×62=2294", $true, $false, $false, $false, $false, $true, 1, $false, "12×81=972", 2) | Out-Null
$d.Content.Find.Execute("17×40=680", $true, $false, $false, $false, $false, $true, 1, $false, "44×35=1540", 2) | Out-Null
$d.Content.Find.Execute("57×86=4902", $true, $false, $false, $false, $false, $true, 1, $false, "50×75=3750", 2) | Out-Null
$d.Content.Find.Execute("20×20=400", $true, $false, $false, $false, $false, $true, 1, $false, "21×20=420", 2) | Out-Null
$d.Content.Find.Execute("87×77=6699", $true, $false, $false, $false, $false, $true, 1, $false, "59×72=4248", 2) | Out-Null
$d.Content.Find.Execute("36×20=720", $true, $false, $false, $false, $false, $true, 1, $false, "72×96=6912", 2) | Out-Null
$d.Content.Find.Execute("25×49=1225", $true, $false, $false, $false, $false, $true, 1, $false, "56×97=5432", 2) | Out-Null
$d.Content.Find.Execute("69×59=4071", $true, $false, $false, $false, $false, $true, 1, $false, "63×40=2520", 2) | Out-Null
$d.Content.Find.Execute("25×75=1875", $true, $false, $false, $false, $false, $true, 1, $false, "36×56=2016", 2) | Out-Null
$d.Content.Find.Execute("57×85=4845", $true, $false, $false, $false, $false, $true, 1, $false, "94×48=4512", 2) | Out-Null
$d.Content.Find.Execute("61×78=4758", $true, $false, $false, $false, $false, $true, 1, $false, "21×92=1932", 2) | Out-Null
$d.Content.Find.Execute("77×12=924", $true, $false, $false, $false, $false, $true, 1, $false, "75×35=2625", 2) | Out-Null
$d.Content.Find.Execute("38×76=2888", $true, $false, $false, $false, $false, $true, 1, $false, "82×88=7216", 2) | Out-Null
